$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update AZ's PopDev value (D2) from 0.0175 to 0.0095
$ws.Range("D2").Value = 0.0095

# Move the active cell selection on the bottom-right pane to D3
$ws.Range("D3").Select()
